# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 392 in the "Naranja" (Orange)
# price sheet, pushing the existing rows 392-407 down to 393-408.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 392 (shifts rows 392:407 -> 393:408)
$ws.Rows.Item(392).Insert()

# Populate the newly inserted row with the new price record
$ws.Range("A392").Value = 2
$ws.Range("B392").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C392").Value = "Coquimbo"
$ws.Range("D392").Value = 44714
$ws.Range("E392").Value = 4
$ws.Range("F392").Value = "Fruta"
$ws.Range("G392").Value = 100102
$ws.Range("H392").Value = "Cítricos"
$ws.Range("I392").Value = 100102005
$ws.Range("J392").Value = "Naranja"
$ws.Range("K392").Value = "Fukumoto"
$ws.Range("L392").Value = "Primera"
$ws.Range("M392").Value = 20
$ws.Range("N392").Value = 180000
$ws.Range("O392").Value = 190000
$ws.Range("P392").Value = 185000
$ws.Range("Q392").Value = "`$/bins (400 kilos)"
$ws.Range("R392").Value = "Provincia de Limarí"
$ws.Range("S392").Value = 462
$ws.Range("T392").Value = 400
